$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.904.15"
$ws.Cells.Item(2, 5).Value = "  +0.99%  "

$ws.Cells.Item(3, 4).Value = "1.628.24"
$ws.Cells.Item(3, 5).Value = "  +1.88%  "

$ws.Cells.Item(4, 5).Value = "  +0.06%  "

$ws.Cells.Item(5, 4).Value = "214.14"
$ws.Cells.Item(5, 5).Value = "  +0.93%  "

$ws.Cells.Item(6, 5).Value = "  +1.04%  "

$ws.Cells.Item(7, 5).Value = "  +0.04%  "

$ws.Cells.Item(8, 5).Value = "  +10.76%  "

$ws.Cells.Item(9, 5).Value = "  +3.31%  "

$ws.Cells.Item(10, 5).Value = "  +2.35%  "

$ws.Cells.Item(11, 4).Value = "0.0915"
$ws.Cells.Item(11, 5).Value = "  +0.85%  "

$ws.Cells.Item(12, 4).Value = "1.861.70"
$ws.Cells.Item(12, 5).Value = "  +2.00%  "

$ws.Cells.Item(13, 4).Value = "1.634.04"
$ws.Cells.Item(13, 5).Value = "  +1.48%  "

$ws.Cells.Item(14, 5).Value = "  +6.15%  "

$ws.Cells.Item(15, 4).Value = "9.23"
$ws.Cells.Item(15, 5).Value = "  +21.58%  "

$ws.Cells.Item(16, 5).Value = "  +3.65%  "

$ws.Cells.Item(17, 4).Value = "29.916.54"
$ws.Cells.Item(17, 5).Value = "  +1.08%  "

$ws.Cells.Item(18, 4).Value = "64.93"
$ws.Cells.Item(18, 5).Value = "  +1.76%  "

$ws.Cells.Item(19, 4).Value = "248.58"
$ws.Cells.Item(19, 5).Value = "  +2.79%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0703"
$ws.Cells.Item(20, 5).Value = "  +1.53%  "

$ws.Cells.Item(21, 5).Value = "  +0.01%  "

$ws.Cells.Item(22, 4).Value = "4.13"
$ws.Cells.Item(22, 5).Value = "  +3.99%  "

$ws.Cells.Item(23, 4).Value = "9.59"
$ws.Cells.Item(23, 5).Value = "  +3.92%  "

$ws.Cells.Item(24, 5).Value = "  +0.66%  "

$ws.Cells.Item(25, 4).Value = "159.05"
$ws.Cells.Item(25, 5).Value = "  +2.92%  "

$ws.Cells.Item(26, 4).Value = "15.69"

$ws.Cells.Item(27, 5).Value = "  +2.19%  "

$ws.Cells.Item(28, 4).Value = "6.58"
$ws.Cells.Item(28, 5).Value = "  +2.98%  "

$ws.Cells.Item(29, 5).Value = "  +0.04%  "

$ws.Cells.Item(30, 4).Value = "0.0489"
$ws.Cells.Item(30, 5).Value = "  +2.51%  "

$ws.Cells.Item(31, 5).Value = "  +5.71%  "

$ws.Cells.Item(32, 5).Value = "  +4.40%  "

$ws.Cells.Item(33, 4).Value = "3.19"
$ws.Cells.Item(33, 5).Value = "  +1.80%  "

$ws.Cells.Item(34, 4).Value = "1.427.87"
$ws.Cells.Item(34, 5).Value = "  -0.25%  "

$ws.Cells.Item(35, 5).Value = "  +7.07%  "

$ws.Cells.Item(36, 4).Value = "1.03"
$ws.Cells.Item(36, 5).Value = "  +1.05%  "

$ws.Cells.Item(37, 5).Value = "  +0.19%  "

$ws.Cells.Item(38, 5).Value = "  -0.12%  "

$ws.Cells.Item(39, 5).Value = "  +3.18%  "

$ws.Cells.Item(40, 5).Value = "  +2.34%  "

$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(41, 4).Value = "71.52"
$ws.Cells.Item(41, 5).Value = "  +8.95%  "

$ws.Cells.Item(42, 2).Value = "BitcoinSV"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(42, 4).Value = "55.23"
$ws.Cells.Item(42, 5).Value = "  +1.64%  "

$ws.Cells.Item(43, 2).Value = "ARBITRUM"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(43, 4).Value = "0.827"
$ws.Cells.Item(43, 5).Value = "  +3.26%  "

$ws.Cells.Item(44, 2).Value = "Kaspa"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(44, 4).Value = "0.0497"
$ws.Cells.Item(44, 5).Value = "  +1.00%  "

$ws.Cells.Item(45, 5).Value = "  +0.54%  "

$ws.Cells.Item(46, 5).Value = "  +5.84%  "

$ws.Cells.Item(47, 5).Value = "  +0.01%  "

$ws.Cells.Item(48, 5).Value = "  +2.39%  "

$ws.Cells.Item(49, 4).Value = "1.768.17"
$ws.Cells.Item(49, 5).Value = "  +1.64%  "

$ws.Cells.Item(50, 4).Value = "89.55"
$ws.Cells.Item(50, 5).Value = "  +4.06%  "

$ws.Cells.Item(51, 5).Value = "  +4.09%  "
